$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as text (e.g. "90.430.48", "3.059.17")
# rather than numbers. For the updated prices that would otherwise look like
# plain decimals (e.g. "243.21"), force the cell to Text format first so Excel
# keeps the new price as a literal string instead of converting it to a number.
$textCells = @('D5','D6','D7','D8','D11','D12','D13','D14','D15','D19','D20','D21','D22','D23','D24','D25','D26','D27','D30','D31','D32','D33','D35','D36','D37','D38','D39','D40','D41','D42','D43','D44','D46','D47','D48','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '90.546.28'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.082.95'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '243.21'
$ws.Range('E5').Value = '  +3.61%  '
$ws.Range('D6').Value = '617.60'
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('D7').Value = '1.13'
$ws.Range('E7').Value = '  +5.15%  '
$ws.Range('D8').Value = '0.364'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '3.080.87'
$ws.Range('E10').Value = '  +6.95%  '
$ws.Range('D11').Value = '0.732'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '0.202'
$ws.Range('E12').Value = '  +2.81%  '
$ws.Range('D13').Value = '0.0000245'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').Value = '34.90'
$ws.Range('E14').Value = '  -5.54%  '
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').Value = '5.47'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '90.433.68'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '3.622.76'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '3.091.16'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').Value = '3.66'
$ws.Range('E19').Value = '  -4.03%  '
$ws.Range('D20').Value = '14.52'
$ws.Range('E20').Value = '  +2.38%  '
$ws.Range('D21').Value = '0.0000210'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '5.75'
$ws.Range('E22').Value = '  +3.37%  '
$ws.Range('D23').Value = '440.30'
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('D24').Value = '9.02'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '91.38'
$ws.Range('E25').Value = '  +2.90%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').Value = '5.60'
$ws.Range('E26').Value = '  -5.53%  '
$ws.Range('D27').Value = '11.86'
$ws.Range('E27').Value = '  -5.92%  '
$ws.Range('D28').Value = '3.223.23'
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '0.178'
$ws.Range('E30').Value = '  +10.62%  '
$ws.Range('D31').Value = '0.238'
$ws.Range('E31').Value = '  +22.04%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '9.17'
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '0.168'
$ws.Range('E33').Value = '  +11.57%  '
$ws.Range('E34').Value = '  +29.51%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = '0.955'
$ws.Range('E35').Value = '  +4.58%  '
$ws.Range('B36').Value = 'MantraDAO'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D36').Value = '4.28'
$ws.Range('E36').Value = '  +25.83%  '
$ws.Range('D37').Value = '26.34'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = '7.58'
$ws.Range('E38').Value = '  +6.14%  '
$ws.Range('D39').Value = '1.90'
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('D40').Value = '487.90'
$ws.Range('E40').Value = '  -5.11%  '
$ws.Range('D41').Value = '3.53'
$ws.Range('E41').Value = '  -6.91%  '
$ws.Range('D42').Value = '1.28'
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('D43').Value = '0.418'
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('D44').Value = '22.18'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '153.99'
$ws.Range('E46').Value = '  +2.49%  '
$ws.Range('D47').Value = '1.89'
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').Value = '0.684'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').Value = '1.34'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '44.24'
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('D51').Value = '4.43'
$ws.Range('E51').Value = '  -0.76%  '
